# Update FlashScore odds/stats data for rows 4, 6, and 7 (Sheet1).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 4 ----
$ws.Range("G4").Value = 1.67
$ws.Range("H4").Value = 3.25
$ws.Range("I4").Value = 6.25
$ws.Range("J4").Value = 2.4
$ws.Range("L4").Value = 7
$ws.Range("M4").Value = 1.14
$ws.Range("N4").Value = 5.5
$ws.Range("W4").Value = 4.5
$ws.Range("X4").Value = 6
$ws.Range("AF4").Value = 126
$ws.Range("AG4").Value = 11
$ws.Range("AI4").Value = 23
$ws.Range("AJ4").Value = 81
$ws.Range("AK4").Value = 67
$ws.Range("AL4").Value = 81
$ws.Range("AR4").Value = 81
$ws.Range("AW4").Value = 7.5
$ws.Range("AX4").Value = 41

# ---- Row 6 ----
$ws.Range("I6").Value = 1.44
$ws.Range("W6").Value = 21
$ws.Range("X6").Value = 34
$ws.Range("Y6").Value = 19
$ws.Range("AD6").Value = 9
$ws.Range("AI6").Value = 9
$ws.Range("AL6").Value = 21
$ws.Range("AZ6").Value = 19

# ---- Row 7 ----
$ws.Range("G7").Value = 2.5
$ws.Range("H7").Value = 3.2
$ws.Range("I7").Value = 2.63
$ws.Range("L7").Value = 3.1
$ws.Range("W7").Value = 10
$ws.Range("AA7").Value = 21
$ws.Range("AB7").Value = 29
$ws.Range("AN7").Value = 4.75
